# Sprint Backlog - FINAL report update
# - Fix the "codigo" -> "código" typo in the two "Procura posicionamento..." tasks
# - Thicken the bottom border of the last table row (row 18) to close the table
# - Update the active selection shown when the workbook is reopened

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# --- Fix accent typo in the shared task-name strings -----------------------
$ws.Range("C4").Value = "Procura posicionamento do código da tarefa 1 (todoList)"
$ws.Range("C5").Value = "Procura posicionamento do código da tarefa 1 (todoList)"
$ws.Range("C6").Value = "Procura posicionamento do código da tarefa 1 (todoList)"
$ws.Range("C7").Value = "Procura posicionamento do código da tarefa 2 (adicionar ficheiros)"
$ws.Range("C8").Value = "Procura posicionamento do código da tarefa 2 (adicionar ficheiros)"

# --- Thicken the bottom border that closes the table on row 18 -------------
$closingRow = $ws.Range("C18:E18")
$closingRow.Borders.Item(9).LineStyle = 1
$closingRow.Borders.Item(9).Weight = -4138

# --- Update the saved view/selection ----------------------------------------
$ws.Activate()
$ws.Range("H8").Select()
